$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Layout and Theme vars")
$ws2 = $wb.Worksheets.Item("Changelog")

# --- Changelog: append 5 new changelog rows (adds new shared strings) ---
$ws2.Cells.Item(15, 1).Value = "thickness-decoration --> textDecorationThickness"
$ws2.Cells.Item(16, 1).Value = "style-decoration --> textDecorationStyle"
$ws2.Cells.Item(17, 1).Value = "line-decoration --> textDecorationLine"
$ws2.Cells.Item(18, 1).Value = "color-decoration --> textDecorationColor"
$ws2.Cells.Item(19, 1).Value = "text-decoration --> textDecoration"

# Rows 17 & 18 pick up the same formatting already used by A2/A3
[void]$ws2.Cells.Item(2, 1).Copy()
[void]$ws2.Cells.Item(17, 1).PasteSpecial(-4122)
[void]$ws2.Cells.Item(18, 1).PasteSpecial(-4122)

# Widen column A on the Changelog sheet to fit the new (longer) text
$ws2.Columns.Item(1).ColumnWidth = 41.6667

# --- Layout and Theme vars: mark the textDecoration-related rows as "done" ---
# (copy the "Good" formatting already used lower in the sheet, e.g. E64)
[void]$ws1.Range("E64").Copy()
foreach ($r in 58..63) {
    [void]$ws1.Cells.Item($r, 5).PasteSpecial(-4122)
}

# --- Selection / active sheet bookkeeping ---
[void]$ws2.Range("A20").Select()
[void]$ws1.Activate()
[void]$ws1.Range("D59").Select()
